# Apply the "stuff at the bottom of the sheets" commit:
#  - add "generic" to J2:J5 (pair_kind column for the practice rows)
#  - add a second small table starting at row 27 with a "stim details" header,
#    a sub-header row, and 8 data rows describing video/audio counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 is an empty shared-string placeholder cell in the source file; explicitly
# re-clear it so the load/save round-trip doesn't leave a stray value behind.
$ws.Range("F1").ClearContents()

# --- J2:J5 get the "generic" pair_kind value, matching rows 6-13's kind column ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- New block starting at row 27 ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
